$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.734.49'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '2.103.90'
$ws.Range('E3').Value = '  +5.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5293'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4357'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08949'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.99'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.83%  '
$ws.Range('E11').Value = '  +2.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = '2.106.24'
$ws.Range('E13').Value = '  +5.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.717'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.758'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001132'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06687'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.294'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('D23').Value = '30.802.87'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.291'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.65%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.352.43'
$ws.Range('E26').Value = '  +5.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.566'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.195'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.75%  '
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.166'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.026'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.533'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +14.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02599'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.579'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.522'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06744'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2269'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6813'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.246'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6405'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.216'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.258'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('E50').Value = '  +3.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.192'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.21%  '
